# Update the "取得日時" (acquired timestamp) column on the "ランサーズ" sheet
# for the existing data rows (2-11) from 06:32:13 to 06:40:58.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newValue = "2025-12-15 06:40:58"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newValue
}
